$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Globo"
$ws.Range("B2").Value = "Inter TV Rural"
$ws.Range("C2").Value = "Agricultura"
$ws.Range("D2").Value = "2025-03-31T19:34"
$ws.Range("E2").Value = "Positivo"
$ws.Range("F2").Value = "Com Nota"
$ws.Range("G2").Value = "Teste"
